$wb = $excel.ActiveWorkbook

# --- Sheet "Prov_5_inc": refresh the OR/CI figures (rows 2-71, columns B-D) ---
$wsInc = $wb.Worksheets.Item("Prov_5_inc")

$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,57,58,59,60,61,62,63,64,65,66,67,68,69,70,71)
$bvals = @(1.06,1.0900000000000001,1.28,0.54,1.37,0.89,1.17,0.69,1.1599999999999999,1.03,1.01,0.68,1.39,0.38,1.08,1.0900000000000001,0.06,0.75,0.45,1.58,1.98,1.1499999999999999,0.95,0.77,0.96,0.03,0.24,1.43,0.95,1.23,0.9,1.1499999999999999,0.99,0.78,0.63,0.6,0.82,0.71,0.99,0.56000000000000005,0.98,0.57999999999999996,1.08,1.19,1.35,0.54,0.94,2.17,1.28,1.1000000000000001,2.4300000000000002,1.55,1,1.65,0.94,0.88,1.1399999999999999,1.01,0.27,0.25,0.59,2.87,0.64,1.63,1.84,2.42,0.98,0.21,0.59,0.17)
$cvals = @(0.92,0.85,0.81,0.19,0.53,0.42,0.54,0.33,0.56000000000000005,0.5,0.98,0.48,1.06,0.22,0.91,0.81,0.03,0.39,0.14000000000000001,0.8,0.94,0.6,0.5,0.4,0.93,0.02,0.18,0.99,0.84,0.98,0.56000000000000005,0.67,0.52,0.39,0.28000000000000003,0.31,0.42,0.36,0.96,0.41,0.77,0.49,0.94,0.93,0.53,0.16,0.51,1.07,0.57999999999999996,0.54,1.23,0.77,0.98,1.1100000000000001,0.74,0.79,0.99,0.78,0.13,0.11,0.28999999999999998,1.35,0.25,0.78,0.89,1.17,0.96,0.14000000000000001,0.45,0.12)
$dvals = @(1.21,1.39,2.0299999999999998,1.38,3.67,1.89,2.5,1.45,2.37,2.12,1.03,0.97,1.82,0.63,1.28,1.47,0.1,1.42,1.32,3.12,4.22,2.23,1.81,1.5,0.99,0.06,0.33,2.0699999999999998,1.07,1.55,1.44,1.97,1.9,1.52,1.43,1.1499999999999999,1.55,1.34,1.01,0.77,1.27,0.69,1.23,1.51,3.43,1.8,1.75,4.47,2.83,2.2599999999999998,4.87,3.14,1.03,2.4700000000000002,1.2,0.99,1.32,1.32,0.54,0.56000000000000005,1.1399999999999999,6.36,1.63,3.54,3.97,5.21,1.01,0.32,0.77,0.23)

for ($i = 0; $i -lt $rows.Count; $i++) {
  $r = $rows[$i]
  $wsInc.Cells.Item($r, 2).Value = $bvals[$i]
  $wsInc.Cells.Item($r, 3).Value = $cvals[$i]
  $wsInc.Cells.Item($r, 4).Value = $dvals[$i]
}

# Group label for the first province block changed from "NFL & NS" to "NL & NS"
$grpRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15)
for ($i = 0; $i -lt $grpRows.Count; $i++) {
  $wsInc.Cells.Item($grpRows[$i], 5).Value = "NL & NS"
}

# Row label fix: "Traveling Distance (km)" -> "Traveling Distance (per 25km)"
$labelRows = @(13,27,41,55,69)
for ($i = 0; $i -lt $labelRows.Count; $i++) {
  $wsInc.Cells.Item($labelRows[$i], 1).Value = "Traveling Distance (per 25km)"
}

# Sheet "Prov_5_inc" becomes the active sheet/tab, with the given selection
$wsInc.Activate()
$wsInc.Range("I58").Select()

# Sheet "Reg" keeps its own cell selection untouched (losing tabSelected happens
# automatically once another sheet is activated above)
$wsReg = $wb.Worksheets.Item("Reg")
$wsReg.Range("C5").Select()
